$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Semana_8")

# Map of row number -> new value for column Q (uds. Objetivo semana pasada)
$updates = @{
    3  = 2
    8  = 2
    9  = 1
    10 = 2
    18 = 2
    19 = 6
    26 = 1
    29 = 2
    32 = 1
    34 = 3
    38 = 1
    41 = 2
    42 = 1
    43 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("Q$row").Value = $updates[$row]
}
